$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 7

# --- Plain numeric cells ---
$ws.Cells.Item($row, 1).Value = 112507258
$ws.Cells.Item($row, 2).Value = 85323
$ws.Cells.Item($row, 5).Value = 1980
$ws.Cells.Item($row, 17).Value = 573889
$ws.Cells.Item($row, 18).Value = 6303251
$ws.Cells.Item($row, 19).Value = 10

# --- Plain text cells ---
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
$ws.Cells.Item($row, 4).Value = "VU"
$ws.Cells.Item($row, 6).Value = "Porslinsblå spindling"
$ws.Cells.Item($row, 7).Value = "Cortinarius cumatilis"
$ws.Cells.Item($row, 8).Value = "Fr."
$ws.Cells.Item($row, 10).Value = "fruktkroppar"
$ws.Cells.Item($row, 16).Value = "Snärjebäcken kalkbarrskog, Sm"
$ws.Cells.Item($row, 20).Value = "Kalmar"
$ws.Cells.Item($row, 21).Value = "Kalmar"
$ws.Cells.Item($row, 22).Value = "Småland"
$ws.Cells.Item($row, 23).Value = "Ryssby"
$ws.Cells.Item($row, 29).Value = "Kalkrik blandskog utmed Snärjebäcken. Blåsippor hassel idegran."
$ws.Cells.Item($row, 49).Value = "Aron Edman"
$ws.Cells.Item($row, 50).Value = "Aron Edman"

# --- Text cells that look like numbers/dates: force literal text via a
# leading quote-prefix (same trick Excel's UI uses), then strip the style
# flag the quote prefix leaves behind so the cell keeps the default style. ---
$ws.Cells.Item($row, 9).Value = "'1"
$ws.Cells.Item($row, 9).Style = "Normal"

$ws.Cells.Item($row, 25).Value = "'2023-10-03"
$ws.Cells.Item($row, 25).Style = "Normal"

$ws.Cells.Item($row, 27).Value = "'2023-10-03"
$ws.Cells.Item($row, 27).Style = "Normal"

# --- Boolean cells ---
$ws.Cells.Item($row, 30).Value = $false
$ws.Cells.Item($row, 31).Value = $false
$ws.Cells.Item($row, 33).Value = $false

# --- Empty text cells (present in the row but with no content) ---
$ws.Cells.Item($row, 11).Value = "'"
$ws.Cells.Item($row, 11).Style = "Normal"

$ws.Cells.Item($row, 46).Value = "'"
$ws.Cells.Item($row, 46).Style = "Normal"

$ws.Cells.Item($row, 51).Value = "'"
$ws.Cells.Item($row, 51).Style = "Normal"
